$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.466326861123889
$ws.Range("D2").Value = 0.03958538122239474
$ws.Range("E2").Value = 0.06367431157244452
$ws.Range("F2").Value = 6.676538951618966
$ws.Range("G2").Value = 0.002680767426274793
$ws.Range("J2").Value = 0.2144512998220023
$ws.Range("K2").Value = 1.334000927136486
$ws.Range("L2").Value = 0.3224214161660939
$ws.Range("B3").Value = 1.461440581020781
$ws.Range("D3").Value = 0.034624188147518
$ws.Range("E3").Value = 0.06344276164392859
$ws.Range("F3").Value = 6.487217307807725
$ws.Range("G3").Value = 0.002686487173836127
$ws.Range("J3").Value = 0.2115868418418749
$ws.Range("K3").Value = 1.313392082407518
$ws.Range("L3").Value = 0.3233204852921148
$ws.Range("B4").Value = 1.459387964095555
$ws.Range("D4").Value = 0.0315711425036227
$ws.Range("E4").Value = 0.06329753118593784
$ws.Range("F4").Value = 6.37171629586993
$ws.Range("G4").Value = 0.002690181294898854
$ws.Range("J4").Value = 0.2097969764911838
$ws.Range("K4").Value = 1.302104696073826
$ws.Range("L4").Value = 0.3241938423872952
$ws.Range("B5").Value = 1.45878988887722
$ws.Range("D5").Value = 0.03032512222974049
$ws.Range("E5").Value = 0.06323756579672768
$ws.Range("F5").Value = 6.324831445438861
$ws.Range("G5").Value = 0.002691732651452601
$ws.Range("J5").Value = 0.2090596529514777
$ws.Range("K5").Value = 1.297847686488268
$ws.Range("L5").Value = 0.3246304288322222
$ws.Range("B6").Value = 1.458704980041261
$ws.Range("D6").Value = 0.03011810345370236
$ws.Range("E6").Value = 0.0632275609081614
$ws.Range("F6").Value = 6.317057169088514
$ws.Range("G6").Value = 0.002691993034184239
$ws.Range("J6").Value = 0.2089367378218938
$ws.Range("K6").Value = 1.297161484278575
$ws.Range("L6").Value = 0.3247077937091234
$ws.Range("B7").Value = 1.459378932901814
$ws.Range("D7").Value = 0.03155434605145047
$ws.Range("E7").Value = 0.06329672565859878
$ws.Range("F7").Value = 6.37108325593104
$ws.Range("G7").Value = 0.002690202030717925
$ws.Range("J7").Value = 0.2097870649932645
$ws.Range("K7").Value = 1.3020458982393
$ws.Range("L7").Value = 0.3241994037994331
$ws.Range("B8").Value = 1.464445455863938
$ws.Range("D8").Value = 0.0378760887763292
$ws.Range("E8").Value = 0.06359509938244234
$ws.Range("F8").Value = 6.611104026415916
$ws.Range("G8").Value = 0.002682701883713598
$ws.Range("J8").Value = 0.2134700090986676
$ws.Range("K8").Value = 1.326610793713627
$ws.Range("L8").Value = 0.3226646532218354
$ws.Range("B9").Value = 1.48189833754131
$ws.Range("D9").Value = 0.05022620431992664
$ws.Range("E9").Value = 0.06415659284751207
$ws.Range("F9").Value = 7.087893104747536
$ws.Range("G9").Value = 0.002669432124654762
$ws.Range("J9").Value = 0.2204515035295813
$ws.Range("K9").Value = 1.385673984544269
$ws.Range("L9").Value = 0.3222105579430661
$ws.Range("B10").Value = 1.49930747183933
$ws.Range("D10").Value = 0.05928255192793586
$ws.Range("E10").Value = 0.06455564948053061
$ws.Range("F10").Value = 7.442243346990466
$ws.Range("G10").Value = 0.002660549019360047
$ws.Range("J10").Value = 0.2254424290246391
$ws.Range("K10").Value = 1.435784594895324
$ws.Range("L10").Value = 0.323444391816011
$ws.Range("B11").Value = 1.508224594396296
$ws.Range("D11").Value = 0.06340124654663271
$ws.Range("E11").Value = 0.06473446916241787
$ws.Range("F11").Value = 7.604400330647252
$ws.Range("G11").Value = 0.002656693712033554
$ws.Range("J11").Value = 0.2276847154870509
$ws.Range("K11").Value = 1.46005741914243
$ws.Range("L11").Value = 0.3243482811191853
$ws.Range("B12").Value = 1.511744776808314
$ws.Range("D12").Value = 0.06496091869037457
$ws.Range("E12").Value = 0.06480180972423399
$ws.Range("F12").Value = 7.66594849309439
$ws.Range("G12").Value = 0.002655260334331739
$ws.Range("J12").Value = 0.228529910956162
$ws.Range("K12").Value = 1.469462657890205
$ws.Range("L12").Value = 0.3247400018866529
$ws.Range("B13").Value = 1.510980261824557
$ws.Range("D13").Value = 0.06462501217649219
$ws.Range("E13").Value = 0.06478732313794655
$ws.Range("F13").Value = 7.652686568630997
$ws.Range("G13").Value = 0.002655567859731978
$ws.Range("J13").Value = 0.2283480547566086
$ws.Range("K13").Value = 1.46742754965598
$ws.Range("L13").Value = 0.3246534362851605
$ws.Range("B14").Value = 1.508511326662131
$ws.Range("D14").Value = 0.06352956050905334
$ws.Range("E14").Value = 0.06474001672771212
$ws.Range("F14").Value = 7.609461045870091
$ws.Range("G14").Value = 0.002656575256334257
$ws.Range("J14").Value = 0.2277543278209677
$ws.Range("K14").Value = 1.460826904886517
$ws.Range("L14").Value = 0.3243795164242727
$ws.Range("B15").Value = 1.507017715897234
$ws.Range("D15").Value = 0.0628585712803158
$ws.Range("E15").Value = 0.06471099187228324
$ws.Range("F15").Value = 7.583002903636952
$ws.Range("G15").Value = 0.002657195767904959
$ws.Range("J15").Value = 0.2273901477158944
$ws.Range("K15").Value = 1.456811679284897
$ws.Range("L15").Value = 0.3242181759822671
$ws.Range("B16").Value = 1.49874479171163
$ws.Range("D16").Value = 0.05901337805003948
$ws.Range("E16").Value = 0.06454390992382697
$ws.Range("F16").Value = 7.431665742052559
$ws.Range("G16").Value = 0.002660804695648875
$ws.Range("J16").Value = 0.2252953359794887
$ws.Range("K16").Value = 1.434228138383673
$ws.Range("L16").Value = 0.3233922291721854
$ws.Range("B17").Value = 1.493925134590143
$ws.Range("D17").Value = 0.05665431794614051
$ws.Range("E17").Value = 0.06444072665748868
$ws.Range("F17").Value = 7.339074756259095
$ws.Range("G17").Value = 0.002663066098240741
$ws.Range("J17").Value = 0.2240031310109174
$ws.Range("K17").Value = 1.420753085231297
$ws.Range("L17").Value = 0.322973407944886
$ws.Range("B18").Value = 1.491246894923222
$ws.Range("D18").Value = 0.05529733689422756
$ws.Range("E18").Value = 0.06438112204682689
$ws.Range("F18").Value = 7.285909071605317
$ws.Range("G18").Value = 0.002664384281722959
$ws.Range("J18").Value = 0.2232572350045672
$ws.Range("K18").Value = 1.413141562474465
$ws.Range("L18").Value = 0.3227647507689682
$ws.Range("B19").Value = 1.490356216390182
$ws.Range("D19").Value = 0.05483786209158836
$ws.Range("E19").Value = 0.06436089645556553
$ws.Range("F19").Value = 7.267923437300226
$ws.Range("G19").Value = 0.002664833603573058
$ws.Range("J19").Value = 0.2230042276060047
$ws.Range("K19").Value = 1.410588263507606
$ws.Range("L19").Value = 0.3226996348013529
$ws.Range("B20").Value = 1.494428477070556
$ws.Range("D20").Value = 0.05690545431075122
$ws.Range("E20").Value = 0.06445173710888774
$ws.Range("F20").Value = 7.348921856128754
$ws.Range("G20").Value = 0.002662823559617734
$ws.Range("J20").Value = 0.2241409621458175
$ws.Range("K20").Value = 1.422173135907741
$ws.Range("L20").Value = 0.3230146544853625
$ws.Range("B21").Value = 1.509232619154204
$ws.Range("D21").Value = 0.06385131955255474
$ws.Range("E21").Value = 0.06475392180746464
$ws.Range("F21").Value = 7.622153516303513
$ws.Range("G21").Value = 0.00265627864027357
$ws.Range("J21").Value = 0.2279288248624951
$ws.Range("K21").Value = 1.462759866052892
$ws.Range("L21").Value = 0.3244586302192118
$ws.Range("B22").Value = 1.519744197646901
$ws.Range("D22").Value = 0.06839099206915478
$ws.Range("E22").Value = 0.0649492411039887
$ws.Range("F22").Value = 7.801560916082906
$ws.Range("G22").Value = 0.002652155800564397
$ws.Range("J22").Value = 0.2303816853432821
$ws.Range("K22").Value = 1.490531528725342
$ws.Range("L22").Value = 0.3256905949321123
$ws.Range("B23").Value = 1.514057451469796
$ws.Range("D23").Value = 0.06596801410347553
$ws.Range("E23").Value = 0.06484518939814432
$ws.Range("F23").Value = 7.705729897148672
$ws.Range("G23").Value = 0.002654342138740046
$ws.Range("J23").Value = 0.2290745829918279
$ws.Range("K23").Value = 1.475594859095338
$ws.Range("L23").Value = 0.3250066400173068
$ws.Range("B24").Value = 1.494200627365501
$ws.Range("D24").Value = 0.05679191783423221
$ws.Range("E24").Value = 0.0644467601669203
$ws.Range("F24").Value = 7.344469776717403
$ws.Range("G24").Value = 0.002662933154997476
$ws.Range("J24").Value = 0.2240786580133829
$ws.Range("K24").Value = 1.421530709145884
$ws.Range("L24").Value = 0.3229959068634827
$ws.Range("B25").Value = 1.476371743253395
$ws.Range("D25").Value = 0.04688902291199781
$ws.Range("E25").Value = 0.06400712755600235
$ws.Range("F25").Value = 6.958221205313237
$ws.Range("G25").Value = 0.00267286908120039
$ws.Range("J25").Value = 0.2185877341734965
$ws.Range("K25").Value = 1.368521958395377
$ws.Range("L25").Value = 0.3220588503189674
